$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Busos
$ws2 = $wb.Worksheets.Item(2)   # Topologia

# --- Sheet1 (Busos): insert a new "delta" column before the old "Tipus" column ---
$ws1.Columns.Item(5).Insert() | Out-Null
$ws1.Range("E1").Value = "delta"
$ws1.Range("E2").Value = 0

# --- Sheet1 (Busos): add new row 14 (bus 12, slack-like row with delta) ---
$ws1.Range("A14").Value = 12
$ws1.Range("D14").Value = 1
$ws1.Range("E14").Value = 0.1
$ws1.Range("F14").Value = "Slack"

# --- Sheet2 (Topologia): add new row 16 ---
$ws2.Range("A16").Value = 5
$ws2.Range("B16").Value = 12
$ws2.Range("C16").Value = 0.01
$ws2.Range("D16").Value = 0.05
$ws2.Range("E16").Value = 0.0001

# --- Selections / active sheet ---
$ws1.Range("F15").Select() | Out-Null
$ws2.Activate()
$ws2.Range("I9").Select() | Out-Null
